$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures scraped for this run.
# Column D cells whose new value is a plain decimal number are forced
# to Text format first so Excel does not silently convert them (the
# sheet stores prices as formatted text, e.g. "62.873.58").

$ws.Range("D2").Value = "62.873.58"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "3.041.86"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.26"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.08"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D9").Value = "3.042.93"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.74"
$ws.Range("E14").Value = "  -4.70%  "
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "3.542.99"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "62.846.62"
$ws.Range("D19").Value = "3.040.80"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.17"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.09"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.54"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("E27").Value = "  +3.86%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.52"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -5.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.84"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("E39").Value = "  -9.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.17"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "422.57"
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").Value = "2.823.30"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0357"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.87"
$ws.Range("E47").Value = "  -5.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.75"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.73"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("E51").Value = "  -1.10%  "
